$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure Price column values stay as text (matches original inlineStr type)
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

$ws.Range('D2').Value = '27.558.35'
$ws.Range('E2').Value = '  -1.33%  '
$ws.Range('D3').Value = '1.843.04'
$ws.Range('E3').Value = '  -2.03%  '
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  -1.19%  '
$ws.Range('D5').Value = '333.61'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('E6').Value = '  -1.09%  '
$ws.Range('D7').Value = '0.4622'
$ws.Range('E7').Value = '  -1.23%  '
$ws.Range('D8').Value = '0.3850'
$ws.Range('E8').Value = '  -1.67%  '
$ws.Range('D9').Value = '45.91'
$ws.Range('E9').Value = '  -2.38%  '
$ws.Range('D10').Value = '0.07894'
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('D11').Value = '0.9943'
$ws.Range('E11').Value = '  -1.52%  '
$ws.Range('D12').Value = '21.48'
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('D13').Value = '1.852.83'
$ws.Range('E13').Value = '  -2.82%  '
$ws.Range('D14').Value = '5.927'
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('D15').Value = '7.119'
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('E16').Value = '  -1.27%  '
$ws.Range('D17').Value = '88.62'
$ws.Range('E17').Value = '  +1.37%  '
$ws.Range('D18').Value = '0.06671'
$ws.Range('E18').Value = '  -1.90%  '
$ws.Range('E19').Value = '  -1.05%  '
$ws.Range('D20').Value = '17.08'
$ws.Range('E20').Value = '  +0.65%  '
$ws.Range('D21').Value = '1.006'
$ws.Range('E21').Value = '  -1.11%  '
$ws.Range('D22').Value = '27.530.49'
$ws.Range('E22').Value = '  -1.49%  '
$ws.Range('D23').Value = '5.373'
$ws.Range('E23').Value = '  -1.69%  '
$ws.Range('E24').Value = '  -0.80%  '
$ws.Range('D25').Value = '2.304'
$ws.Range('E25').Value = '  -2.19%  '
$ws.Range('D26').Value = '2.074.90'
$ws.Range('E26').Value = '  -2.64%  '
$ws.Range('D27').Value = '158.59'
$ws.Range('E27').Value = '  -0.53%  '
$ws.Range('D28').Value = '19.46'
$ws.Range('E28').Value = '  -2.65%  '
$ws.Range('D29').Value = '2.102'
$ws.Range('E29').Value = '  +1.67%  '
$ws.Range('D30').Value = '5.392'
$ws.Range('E30').Value = '  -0.96%  '
$ws.Range('D31').Value = '119.82'
$ws.Range('E31').Value = '  -0.76%  '
$ws.Range('D32').Value = '0.9739'
$ws.Range('E32').Value = '  +1.90%  '
$ws.Range('D33').Value = '0.09375'
$ws.Range('E33').Value = '  -1.61%  '
$ws.Range('E34').Value = '  -1.67%  '
$ws.Range('D35').Value = '5.270'
$ws.Range('E35').Value = '  -0.76%  '
$ws.Range('D36').Value = '1.339'
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('D37').Value = '0.06025'
$ws.Range('E37').Value = '  -1.46%  '
$ws.Range('E38').Value = '  -0.49%  '
$ws.Range('D39').Value = '8.271'
$ws.Range('E39').Value = '  +1.86%  '
$ws.Range('D40').Value = '1.182'
$ws.Range('E40').Value = '  -1.93%  '
$ws.Range('D41').Value = '0.5869'
$ws.Range('E41').Value = '  +0.09%  '
$ws.Range('D42').Value = '0.1863'
$ws.Range('E42').Value = '  -1.51%  '
$ws.Range('D43').Value = '10.26'
$ws.Range('E43').Value = '  +1.09%  '
$ws.Range('D44').Value = '1.238'
$ws.Range('E44').Value = '  -2.68%  '
$ws.Range('D45').Value = '0.5577'
$ws.Range('E45').Value = '  -0.80%  '
$ws.Range('D46').Value = '12.12'
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('E47').Value = '  -0.79%  '
$ws.Range('D48').Value = '0.06698'
$ws.Range('E48').Value = '  -2.38%  '
$ws.Range('D49').Value = '110.89'
$ws.Range('E49').Value = '  -2.13%  '
$ws.Range('D50').Value = '1.051'
$ws.Range('E50').Value = '  -1.28%  '
$ws.Range('E51').Value = '  -1.18%  '
